$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.864.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.21%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.651.62'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.19%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.41'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.60%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.27%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.19%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.698'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.59%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.49%  '

# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.02'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.00%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000272'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -6.03%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.240.32'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.61%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.653.94'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.86%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.86%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.86%  '

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.97%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.621.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.34%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.42'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.81%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '399.89'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.66%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '87.47'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.77%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.91%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.94'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.44%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.43'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.40%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.64'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.26'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.77'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.94%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.43%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.27'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.30%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '66.22'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.51%  '

# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '44.33'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.39%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '603.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.05%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.392'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.00%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.08%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0770'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -13.48%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.135'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.20%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.03%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.81%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -8.38%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.88%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.787.64'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.40%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.24%  '

# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.18'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.27%  '

# Row 49
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.76'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.90%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.27%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.51'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -14.72%  '
